# Update "Sprint Backlog Burndown.xlsx" to reflect hours worked.
# Target sheet is "Sprint 2" (the active / tab-selected sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Initial Estimate" column (D) for the rows whose remaining
# hours changed (reflecting hours worked during the sprint).
$ws.Range("D6").Value  = 8
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("D19").Value = 3
$ws.Range("D20").Value = 3

# Recalculate so the SUM formula in D30 (and the dependent chart series
# cache) picks up the new totals.
$excel.CalculateFullRebuild()

# Move/update the active selection on the sheet.
$ws.Activate()
$ws.Range("D42").Select()
